$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" ---
# (appears in the "Status" column of every sheet: Overview!E2:F2, zh-cn!C2, de-de!C2)
$wb.Worksheets.Item("Overview").Range("E2").Value = "In Translation"
$wb.Worksheets.Item("Overview").Range("F2").Value = "In Translation"
$wb.Worksheets.Item("zh-cn").Range("C2").Value = "In Translation"
$wb.Worksheets.Item("de-de").Range("C2").Value = "In Translation"

# --- Narrow the "Status" column now that its text is shorter ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Columns.Item(3).ColumnWidth = 12.5

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Columns.Item(3).ColumnWidth = 12.5
